# تعديل يدوي في شيت Card20 by admin at 2025-12-16 11:57:36
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")
$ws.Activate()

# Row 8 updates
$ws.Range("D8").Value = "'938"
$ws.Range("F8").Value = "✅"
$ws.Range("K8").Value = "✅"
$ws.Range("L8").Value = "16/12/2025"

# Row 35 updates
$ws.Range("A35").Value = "nan"
$ws.Range("L35").Value = "nan"
$ws.Range("M35").Value = "nan"
$ws.Range("N35").Value = "nan"
